# Daily attendance processing - 2025-12-31 22:57:30
# Swap the order of names in the "Recorded By" column (G): cells that read
# "System, dnasr281@gmail.com" become "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value()
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
